# "Generate Report for Handback": refresh the handback timestamps for the
# de0447c2... file (row 3) on both locale sheets now that it has been
# handed back again.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D3").Value = "2016-01-27 08:16:36"   # Correspond Handoff Datetime
$wsZh.Range("G3").Value = "2016-01-27 08:17:23"   # Correspond Handback DateTime

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D3").Value = "2016-01-27 08:16:48"   # Correspond Handoff Datetime
$wsDe.Range("G3").Value = "2016-01-27 08:17:44"   # Correspond Handback DateTime
